$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then optional B/C/D/E new text values.
# 'DForce' marks D-values that parse as plain numbers and must be
# pinned to text (NumberFormat "@") before assignment so Excel's
# COM type-sniffing doesn't silently convert them to Double.
$updates = @(
    @{ Row = 2; D = '42.937.09'; DForce = $false; E = '  -0.59%  ' }
    @{ Row = 3; D = '2.559.43'; DForce = $false; E = '  +0.39%  ' }
    @{ Row = 4; E = '  +0.15%  ' }
    @{ Row = 5; D = '315.71'; DForce = $true; E = '  -0.91%  ' }
    @{ Row = 6; D = '96.60'; DForce = $true; E = '  +0.04%  ' }
    @{ Row = 7; D = '0.576'; DForce = $true; E = '  -0.89%  ' }
    @{ Row = 8; E = '  +0.04%  ' }
    @{ Row = 9; D = '0.538'; DForce = $true; E = '  +0.63%  ' }
    @{ Row = 10; D = '35.57'; DForce = $true; E = '  -3.00%  ' }
    @{ Row = 11; D = '0.0813'; DForce = $true; E = '  -0.26%  ' }
    @{ Row = 12; D = '7.44'; DForce = $true; E = '  -3.02%  ' }
    @{ Row = 13; E = '  -4.64%  ' }
    @{ Row = 14; D = '2.954.62'; DForce = $false; E = '  +0.64%  ' }
    @{ Row = 15; D = '2.509.88'; DForce = $false; E = '  -1.63%  ' }
    @{ Row = 16; D = '15.09'; DForce = $true; E = '  -3.61%  ' }
    @{ Row = 17; D = '0.842'; DForce = $true; E = '  -1.66%  ' }
    @{ Row = 18; D = '42.984.11'; DForce = $false; E = '  +0.00%  ' }
    @{ Row = 19; D = '6.82'; DForce = $true; E = '  +2.74%  ' }
    @{ Row = 20; D = '12.54'; DForce = $true; E = '  -4.26%  ' }
    @{ Row = 21; D = '0.0₃0959'; DForce = $false; E = '  -1.37%  ' }
    @{ Row = 22; D = '69.34'; DForce = $true; E = '  -1.78%  ' }
    @{ Row = 23; D = '252.30'; DForce = $true; E = '  -0.19%  ' }
    @{ Row = 24; E = '  -0.39%  ' }
    @{ Row = 25; D = '2.06'; DForce = $true; E = '  +1.56%  ' }
    @{ Row = 26; D = '26.72'; DForce = $true; E = '  -1.16%  ' }
    @{ Row = 27; E = '  -0.05%  ' }
    @{ Row = 28; E = '  +0.66%  ' }
    @{ Row = 29; D = '39.87'; DForce = $true; E = '  -0.36%  ' }
    @{ Row = 30; D = '10.17'; DForce = $true; E = '  -0.84%  ' }
    @{ Row = 31; D = '5.81'; DForce = $true; E = '  -4.82%  ' }
    @{ Row = 32; D = '154.70'; DForce = $true; E = '  -0.48%  ' }
    @{ Row = 33; E = '  +1.68%  ' }
    @{ Row = 34; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.0804'; DForce = $true; E = '  +1.33%  ' }
    @{ Row = 35; B = 'WEMIXToken'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '2.71'; DForce = $true; E = '  +2.69%  ' }
    @{ Row = 36; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '2.12'; DForce = $true; E = '  -1.66%  ' }
    @{ Row = 37; D = '19.07'; DForce = $true; E = '  -0.55%  ' }
    @{ Row = 38; E = '  -0.87%  ' }
    @{ Row = 39; E = '  +5.75%  ' }
    @{ Row = 40; E = '  -1.05%  ' }
    @{ Row = 41; D = '22.53'; DForce = $true; E = '  -5.68%  ' }
    @{ Row = 42; D = '3.94'; DForce = $true; E = '  +2.49%  ' }
    @{ Row = 43; E = '  -0.21%  ' }
    @{ Row = 44; E = '  +0.24%  ' }
    @{ Row = 45; D = '3.26'; DForce = $true; E = '  -3.12%  ' }
    @{ Row = 46; D = '2.001.63'; DForce = $false; E = '  -1.14%  ' }
    @{ Row = 47; D = '9.01'; DForce = $true; E = '  +1.47%  ' }
    @{ Row = 48; D = '83.19'; DForce = $true; E = '  -3.50%  ' }
    @{ Row = 49; D = '2.809.87'; DForce = $false; E = '  +1.00%  ' }
    @{ Row = 50; D = '74.01'; DForce = $true; E = '  -1.06%  ' }
    @{ Row = 51; D = '103.74'; DForce = $true; E = '  +0.50%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value2 = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value2 = $u.C }
    if ($u.ContainsKey('D')) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.DForce) {
            # Pin to text so COM's type-sniffing doesn't coerce a plain
            # decimal-looking string into a Double, then restore the
            # original (unstyled) cell style so only the value changes.
            $dCell.NumberFormat = '@'
            $dCell.Value2 = $u.D
            $dCell.Style = 'Normal'
        } else {
            $dCell.Value2 = $u.D
        }
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value2 = $u.E }
}
